$wb = $excel.ActiveWorkbook

# --- Sheet "python": refresh the puzzle grid with new values ---
$ws = $wb.Worksheets.Item("python")
$ws.Range("A2").Value = 9
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 7
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()

$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 5
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = 9
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = 1
$ws.Range("I3").ClearContents()

$ws.Range("A4").Value = 7
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 5
$ws.Range("G4").ClearContents()
$ws.Range("H4").Value = 9
$ws.Range("I4").Value = 8

$ws.Range("A5").Value = 5
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()

$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = 1
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 6
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = 8
$ws.Range("I6").ClearContents()

$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 1
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = 9

$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 5
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 2
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = 4

$ws.Range("A9").ClearContents()
$ws.Range("B9").Value = 6
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").Value = 8
$ws.Range("F9").ClearContents()
$ws.Range("G9").Value = 9
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()

$ws.Range("A10").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 6
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("I10").Value = 2

# New active cell / selection on the "python" sheet
$ws.Activate()
$ws.Range("I12").Select()
$excel.ActiveWindow.Zoom = 88

# --- Remaining sheets: only the zoom level changes (168 -> 88); selection stays put ---
$zoomSheets = @("SDK_1", "SDK_2", "SDK_3", "SDK_4", "Trip square", "NRC", "Trucs")
foreach ($name in $zoomSheets) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Activate()
    $excel.ActiveWindow.Zoom = 88
}

# Restore "python" as the active/selected sheet (it was tabSelected in the original file)
$ws.Activate()
